# Add season-record columns (Wins / Losses / Ties) to the BAL_2018 sheet.
# The source data previously only captured team/player statistics and never
# downloaded the season record, so this backfills the team's 2018 season
# record (47 wins, 115 losses, 0 ties) across every row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (copy formatting from the neighboring header cell so the new
# headers match the existing bold/bordered/centered header style)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-58: same season record repeated for every player on the roster
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 47
    $ws.Cells.Item($r, 31).Value = 115
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-58"
